# "Updates to Ranks, Update LF pathway, etc."
#
# The worksheet is a flat data table (ReachName, Basin, ... HQ_Score_Protection)
# that grew from 11 data rows (A2:W12) to 15 data rows (A2:W16), plus the
# P/Q header columns ("Riparian-Disturbance_score" / "Riparian-CanopyCover_score")
# were swapped, several existing rows got re-scored, and several reaches were
# renamed / re-ordered (Entiat/Methow/Nason/Twisp reaches shuffled and two new
# "Nason Creek Lower" rows plus three new "Twisp River Middle" rows were added).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap the two Riparian columns ---
$ws.Range("P1").Value = "Riparian-CanopyCover_score"
$ws.Range("Q1").Value = "Riparian-Disturbance_score"

# --- Row 2: Big Meadow Creek 01 (P/Q swap carries the old P2 value into Q2) ---
$ws.Range("P2").Value = $null
$ws.Range("Q2").Value = 1

# --- Row 3: Entiat River Potato 05 ---
$ws.Range("N3").Value = 5
$ws.Range("T3").Value = 39
$ws.Range("U3").Value = 0.8666666666666667
$ws.Range("V3").Value = 3

# --- Row 4: Entiat River Potato 06 ---
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 5

# --- Row 5: was "Methow River Fawn 04" -> now "Entiat River Potato 07" ---
$ws.Range("A5").Value = "Entiat River Potato 07"
$ws.Range("B5").Value = "Entiat"
$ws.Range("C5").Value = "Entiat River-Potato Creek"
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 3
$ws.Range("M5").Value = 3
$ws.Range("N5").Value = 3
$ws.Range("S5").Value = 3

# --- Row 6: was "Methow River Thompson 08" -> now "Entiat River Potato 08" ---
$ws.Range("A6").Value = "Entiat River Potato 08"
$ws.Range("B6").Value = "Entiat"
$ws.Range("C6").Value = "Entiat River-Potato Creek"
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 3
$ws.Range("L6").Value = 3
$ws.Range("P6").Value = 3
$ws.Range("Q6").Value = 3
$ws.Range("R6").Value = 3
$ws.Range("S6").Value = 3
$ws.Range("T6").Value = 33
$ws.Range("U6").Value = 0.7333333333333333

# --- Row 7: was "Nason Creek Lower 01" -> now "Methow River Fawn 04" ---
$ws.Range("A7").Value = "Methow River Fawn 04"
$ws.Range("B7").Value = "Methow"
$ws.Range("C7").Value = "Methow River-Fawn Creek"
$ws.Range("G7").Value = 5
$ws.Range("I7").Value = 5
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 3
$ws.Range("M7").Value = 5
$ws.Range("Q7").Value = 5
$ws.Range("R7").Value = 4

# --- Row 8: was "Nason Creek Lower 02" -> now "Methow River Thompson 07" ---
$ws.Range("A8").Value = "Methow River Thompson 07"
$ws.Range("B8").Value = "Methow"
$ws.Range("C8").Value = "Methow River-Thompson Creek"
$ws.Range("L8").Value = 1
$ws.Range("N8").Value = 3
$ws.Range("Q8").Value = 3
$ws.Range("R8").Value = 4
$ws.Range("S8").Value = 5
$ws.Range("T8").Value = 35
$ws.Range("U8").Value = 0.7777777777777778
$ws.Range("V8").Value = 5

# --- Row 9: was "Nason Creek Lower 03" -> now "Methow River Thompson 08" ---
$ws.Range("A9").Value = "Methow River Thompson 08"
$ws.Range("B9").Value = "Methow"
$ws.Range("C9").Value = "Methow River-Thompson Creek"
$ws.Range("L9").Value = 1
$ws.Range("N9").Value = 3
$ws.Range("P9").Value = 5
$ws.Range("R9").Value = 4
$ws.Range("S9").Value = 5
$ws.Range("T9").Value = 35
$ws.Range("U9").Value = 0.7777777777777778

# --- Row 10: was "Twisp River Middle 01" -> now "Methow River Thompson 09" ---
$ws.Range("A10").Value = "Methow River Thompson 09"
$ws.Range("C10").Value = "Methow River-Thompson Creek"
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 4
$ws.Range("L10").Value = 1
$ws.Range("N10").Value = 3
$ws.Range("P10").Value = 5
$ws.Range("R10").Value = 4
$ws.Range("S10").Value = 3

# --- Row 11: was "Twisp River Middle 02" -> now "Nason Creek Lower 01" ---
$ws.Range("A11").Value = "Nason Creek Lower 01"
$ws.Range("B11").Value = "Wenatchee"
$ws.Range("C11").Value = "Lower Nason Creek"
$ws.Range("H11").Value = 5
$ws.Range("I11").Value = 4
$ws.Range("P11").Value = 5
$ws.Range("Q11").Value = 5
$ws.Range("R11").Value = 5
$ws.Range("T11").Value = 38
$ws.Range("U11").Value = 0.8444444444444444
$ws.Range("V11").Value = 3

# --- Row 12: was "Twisp River Middle 06" -> now "Nason Creek Lower 02" ---
$ws.Range("A12").Value = "Nason Creek Lower 02"
$ws.Range("B12").Value = "Wenatchee"
$ws.Range("C12").Value = "Lower Nason Creek"
$ws.Range("G12").Value = 3
$ws.Range("I12").Value = 4
$ws.Range("K12").Value = 5
$ws.Range("M12").Value = 3
$ws.Range("N12").Value = 5
$ws.Range("Q12").Value = 5
$ws.Range("R12").Value = 5
$ws.Range("T12").Value = 40
$ws.Range("U12").Value = 0.8888888888888888
$ws.Range("V12").Value = 3

# --- Row 13: NEW - "Nason Creek Lower 03" ---
$ws.Range("A13").Value = "Nason Creek Lower 03"
$ws.Range("B13").Value = "Wenatchee"
$ws.Range("C13").Value = "Lower Nason Creek"
$ws.Range("D13").Value = "yes"
$ws.Range("E13").Value = "yes"
$ws.Range("F13").Value = "yes"
$ws.Range("G13").Value = 3
$ws.Range("H13").Value = 5
$ws.Range("I13").Value = 4
$ws.Range("J13").Value = 5
$ws.Range("K13").Value = 5
$ws.Range("L13").Value = 3
$ws.Range("M13").Value = 3
$ws.Range("N13").Value = 5
$ws.Range("O13").Value = 5
$ws.Range("P13").Value = 5
$ws.Range("Q13").Value = 5
$ws.Range("R13").Value = 5
$ws.Range("S13").Value = 1
$ws.Range("T13").Value = 36
$ws.Range("U13").Value = 0.8
$ws.Range("V13").Value = 3
$ws.Range("W13").Value = 3

# --- Row 14: NEW - "Twisp River Middle 01" ---
$ws.Range("A14").Value = "Twisp River Middle 01"
$ws.Range("B14").Value = "Methow"
$ws.Range("C14").Value = "Middle Twisp River"
$ws.Range("D14").Value = "yes"
$ws.Range("E14").Value = "yes"
$ws.Range("F14").Value = "yes"
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = 3
$ws.Range("I14").Value = 3
$ws.Range("J14").Value = 5
$ws.Range("K14").Value = 5
$ws.Range("L14").Value = 3
$ws.Range("M14").Value = 3
$ws.Range("N14").Value = 5
$ws.Range("O14").Value = 5
$ws.Range("P14").Value = 3
$ws.Range("Q14").Value = 3
$ws.Range("R14").Value = 3
$ws.Range("S14").Value = 1
$ws.Range("T14").Value = 33
$ws.Range("U14").Value = 0.7333333333333333
$ws.Range("V14").Value = 5
$ws.Range("W14").Value = 3

# --- Row 15: NEW - "Twisp River Middle 02" ---
$ws.Range("A15").Value = "Twisp River Middle 02"
$ws.Range("B15").Value = "Methow"
$ws.Range("C15").Value = "Middle Twisp River"
$ws.Range("D15").Value = "yes"
$ws.Range("E15").Value = "yes"
$ws.Range("F15").Value = "yes"
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 3
$ws.Range("I15").Value = 3
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = 5
$ws.Range("L15").Value = 5
$ws.Range("M15").Value = 3
$ws.Range("N15").Value = 5
$ws.Range("O15").Value = 5
$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 3
$ws.Range("R15").Value = 3
$ws.Range("S15").Value = 1
$ws.Range("T15").Value = 35
$ws.Range("U15").Value = 0.7777777777777778
$ws.Range("V15").Value = 5
$ws.Range("W15").Value = 3

# --- Row 16: NEW - "Twisp River Middle 06" (re-added at the end, re-scored) ---
$ws.Range("A16").Value = "Twisp River Middle 06"
$ws.Range("B16").Value = "Methow"
$ws.Range("C16").Value = "Middle Twisp River"
$ws.Range("D16").Value = "yes"
$ws.Range("E16").Value = "yes"
$ws.Range("F16").Value = "yes"
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 5
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 5
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 5
$ws.Range("M16").Value = 5
$ws.Range("N16").Value = 5
$ws.Range("O16").Value = 5
$ws.Range("P16").Value = 1
$ws.Range("Q16").Value = 5
$ws.Range("R16").Value = 3
$ws.Range("S16").Value = 3
$ws.Range("T16").Value = 37
$ws.Range("U16").Value = 0.8222222222222222
$ws.Range("V16").Value = 3
$ws.Range("W16").Value = 3
